{"js": "// Update the \"Profile\" bullet list describing the REST protocol/endpoints.\n// The document has a flat list of bullet paragraphs (numId=1). Starting at\n// the paragraph \"Protocol: SIDs URNs Resources. ...\" eight paragraphs are\n// rewritten/expanded into eleven paragraphs, ending again with\n// \"Sample Workflow: ToDo.\" followed by the trailing (unchanged) empty\n// paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its (unique) original text so the script\n// is resilient to any unrelated paragraphs before it.\nconst anchorText =\n  \"Protocol: SIDs URNs Resources. Endpoints: Case Classes Events Signatures, Statement Data Pattern Matching Events.\";\nlet idx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    idx = i;\n    break;\n  }\n}\nif (idx === -1) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\nconst items = paragraphs.items;\n\n// 1) \"Protocol: SIDs URNs Resources. ...\" -> rewritten text.\nitems[idx].insertText(\n  \"Protocol: SIDs URNs Resources. Endpoints: Case Classes Aggregated Message Signatures, Aligned Statements Data Pattern Matching Message Events Resource Statement Occurrences. Resource Monad.\",\n  \"Replace\"\n);\n\n// 2) \"Core Model Upper Resources (...)\" -> appended sentence.\nitems[idx + 1].insertText(\n  \"Core Model Upper Resources (DCI Context / Facets: Metaclass, Class, etc. as Resource, root navigation Context Resource). Aggregation (schema cases) / Alignment (resource statements occurrences): Activation.\",\n  \"Replace\"\n);\n\n// 3) NEW paragraph inserted right after (2): \"Core Model Functional Transforms: ...\"\nitems[idx + 1].insertParagraph(\n  \"Core Model Functional Transforms: Functional Activation Statements:  Aggregation Schema Case Classes Statements / Alignment Message Events Resource Statement Occurrence.\",\n  \"After\"\n);\n\n// 4) \"Protocol: GET URN Case Classes / ...\" -> rewritten text.\nitems[idx + 2].insertText(\n  \"Protocol: GET URN Case Classes (Aggregation) / Statements Data (Alignment) Message Events Resource Statement Occurrences.\",\n  \"Replace\"\n);\n\n// 5) \"Protocol: Browse Messages Events Statements. ...\" -> rewritten text.\nitems[idx + 3].insertText(\n  \"Protocol: GET Browse Resource Aggregated / Aligned Message Events Resource Statement Occurrences. Build Context State Flows (Monad Functional Activation).\",\n  \"Replace\"\n);\n\n// 6) \"Protocol: POST URN Navigation Context built ...\" -> rewritten text.\nitems[idx + 4].insertText(\n  \"Protocol: POST URN Navigation Context State Built Resource Activation Data Statements.\",\n  \"Replace\"\n);\n\n// 7) \"Protocol: POST Subsequent entailed Context Browsing / Events Transforms.\" -> rewritten text.\nitems[idx + 5].insertText(\n  \"Protocol: POST Subsequent entailed Context Browsing / Events Functional Transforms Activations.\",\n  \"Replace\"\n);\n\n// 8) & 9) Two brand-new paragraphs inserted after (7), before \"Sample Workflow: ToDo\".\nconst p8 = items[idx + 5].insertParagraph(\n  \"Monad: Resources (Metaclass, Class, etc.). Context.\",\n  \"After\"\n);\np8.insertParagraph(\n  \"Transform: Statements (schema and occurrences).\",\n  \"After\"\n);\n\n// 10) \"Sample Workflow: ToDo\" -> add trailing period. The paragraph has a\n// trailing empty run after the text run; \"Replace\" on the paragraph only\n// rewrites the first (text-bearing) run, leaving the empty run untouched.\nitems[idx + 6].insertText(\"Sample Workflow: ToDo.\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Update the \"Profile\" bullet list describing the REST protocol/endpoints.\n# The document has a flat list of bullet paragraphs (numId=1). Starting at\n# the paragraph \"Protocol: SIDs URNs Resources. ...\" eight paragraphs are\n# rewritten/expanded into eleven paragraphs, ending again with\n# \"Sample Workflow: ToDo.\" followed by the trailing (unchanged) empty\n# paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its (unique) original text so the script\n# is resilient to any unrelated paragraphs before it.\n$anchorText = \"Protocol: SIDs URNs Resources. Endpoints: Case Classes Events Signatures, Statement Data Pattern Matching Events.\"\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  if ($d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\") -eq $anchorText) {\n    $anchorIndex = $i\n    break\n  }\n}\nif ($anchorIndex -eq -1) {\n  throw \"Anchor paragraph not found\"\n}\n\n# 1) \"Protocol: SIDs URNs Resources. ...\" -> rewritten text.\n$p1 = $d.Paragraphs.Item($anchorIndex)\n$p1.Range.Text = \"Protocol: SIDs URNs Resources. Endpoints: Case Classes Aggregated Message Signatures, Aligned Statements Data Pattern Matching Message Events Resource Statement Occurrences. Resource Monad.\"\n\n# 2) \"Core Model Upper Resources (...)\" -> appended sentence.\n$p2 = $d.Paragraphs.Item($anchorIndex + 1)\n$p2.Range.Text = \"Core Model Upper Resources (DCI Context / Facets: Metaclass, Class, etc. as Resource, root navigation Context Resource). Aggregation (schema cases) / Alignment (resource statements occurrences): Activation.\"\n\n# 3) NEW paragraph inserted right after (2): \"Core Model Functional Transforms: ...\"\n$p2.Range.InsertParagraphAfter()\n$p3 = $p2.Next()\n$p3.Range.Text = \"Core Model Functional Transforms: Functional Activation Statements:  Aggregation Schema Case Classes Statements / Alignment Message Events Resource Statement Occurrence.\"\n\n# 4) \"Protocol: GET URN Case Classes / ...\" -> rewritten text.\n$p4 = $d.Paragraphs.Item($anchorIndex + 3)\n$p4.Range.Text = \"Protocol: GET URN Case Classes (Aggregation) / Statements Data (Alignment) Message Events Resource Statement Occurrences.\"\n\n# 5) \"Protocol: Browse Messages Events Statements. ...\" -> rewritten text.\n$p5 = $d.Paragraphs.Item($anchorIndex + 4)\n$p5.Range.Text = \"Protocol: GET Browse Resource Aggregated / Aligned Message Events Resource Statement Occurrences. Build Context State Flows (Monad Functional Activation).\"\n\n# 6) \"Protocol: POST URN Navigation Context built ...\" -> rewritten text.\n$p6 = $d.Paragraphs.Item($anchorIndex + 5)\n$p6.Range.Text = \"Protocol: POST URN Navigation Context State Built Resource Activation Data Statements.\"\n\n# 7) \"Protocol: POST Subsequent entailed Context Browsing / Events Transforms.\" -> rewritten text.\n$p7 = $d.Paragraphs.Item($anchorIndex + 6)\n$p7.Range.Text = \"Protocol: POST Subsequent entailed Context Browsing / Events Functional Transforms Activations.\"\n\n# 8) & 9) Two brand-new paragraphs inserted after (7), before \"Sample Workflow: ToDo\".\n$p7.Range.InsertParagraphAfter()\n$p8 = $p7.Next()\n$p8.Range.Text = \"Monad: Resources (Metaclass, Class, etc.). Context.\"\n\n$p8.Range.InsertParagraphAfter()\n$p9 = $p8.Next()\n$p9.Range.Text = \"Transform: Statements (schema and occurrences).\"\n\n# 10) \"Sample Workflow: ToDo\" -> add trailing period. The paragraph has a\n# trailing empty run after the text run; assigning Range.Text only rewrites\n# the text content, leaving the empty run after it untouched.\n$p10 = $d.Paragraphs.Item($anchorIndex + 9)\n$p10.Range.Text = \"Sample Workflow: ToDo.\"\n"}
